$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1288.375
$ws.Range("I55").Value = 265.5
$ws.Range("J55").Value = 1629.3334
$ws.Range("K55").Value = 265.5
$ws.Range("L55").Value = 1629.3334
$ws.Range("M55").Value = -51.5
$ws.Range("N55").Value = -2057.3334
$ws.Range("H62").Value = 15628158
$ws.Range("I62").Value = 20835944
$ws.Range("K62").Value = 20835944
$ws.Range("M62").Value = -20835320
$ws.Range("H65").Value = 15628158
$ws.Range("I65").Value = 20835944
$ws.Range("K65").Value = 104179720
$ws.Range("M65").Value = -104176600
$ws.Range("H69").Value = 11293
$ws.Range("J69").Value = 13999.667
$ws.Range("L69").Value = 41999.001
$ws.Range("N69").Value = -43747.001
$ws.Range("H72").Value = 11293
$ws.Range("J72").Value = 13999.667
$ws.Range("L72").Value = 125997.003
$ws.Range("N72").Value = -134733.003
$ws.Range("H96").Value = 685.6
$ws.Range("I96").Value = 685.6
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2056.8
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -683.8000000000002
$ws.Range("N96").ClearContents()
$ws.Range("H98").Value = 1775.4849
$ws.Range("I98").Value = 888.11536
$ws.Range("K98").Value = 888.11536
$ws.Range("M98").Value = 609.88464
$ws.Range("H122").Value = 1775.4849
$ws.Range("I122").Value = 888.11536
$ws.Range("K122").Value = 2664.34608
$ws.Range("M122").Value = -214.3460800000003
$ws.Range("H127").Value = 20137.334
$ws.Range("I127").Value = 21579.5
$ws.Range("K127").Value = 64738.5
$ws.Range("M127").Value = -59778.5
$ws.Range("H132").Value = 2099.5144
$ws.Range("I132").Value = 1349.4333
$ws.Range("J132").Value = 6600
$ws.Range("K132").Value = 4048.2999
$ws.Range("L132").Value = 19800
$ws.Range("M132").Value = -1518.2999
$ws.Range("N132").Value = -24860

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7053.7334
$ws.Range("I45").Value = 2923.111
$ws.Range("K45").Value = 2923.111
$ws.Range("M45").Value = -2546.111
$ws.Range("H96").Value = 24672
$ws.Range("J96").Value = 24672
$ws.Range("L96").Value = 24672
$ws.Range("N96").Value = -30164
$ws.Range("H129").Value = 78401
$ws.Range("J129").Value = 78401
$ws.Range("L129").Value = 78401
$ws.Range("N129").Value = -88401
$ws.Range("H132").Value = 3781.2898
$ws.Range("J132").Value = 9462.362999999999
$ws.Range("L132").Value = 28387.089
$ws.Range("N132").Value = -33447.089

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 306.25
$ws.Range("I2").Value = 306.25
$ws.Range("K2").Value = 306.25
$ws.Range("M2").Value = -193.25
$ws.Range("H99").Value = 5045.0557
$ws.Range("J99").Value = 6402.3
$ws.Range("L99").Value = 6402.3
$ws.Range("N99").Value = -9398.299999999999
$ws.Range("H105").Value = 1352.6
$ws.Range("I105").Value = 1107.5
$ws.Range("K105").Value = 1107.5
$ws.Range("M105").Value = 639.5
$ws.Range("H122").Value = 1901.9412
$ws.Range("J122").Value = 3037.6365
$ws.Range("L122").Value = 9112.9095
$ws.Range("N122").Value = -14012.9095
$ws.Range("H126").Value = 5045.0557
$ws.Range("J126").Value = 6402.3
$ws.Range("L126").Value = 19206.9
$ws.Range("N126").Value = -24146.9
$ws.Range("H132").Value = 4003.8276
$ws.Range("I132").Value = 2584.6316
$ws.Range("K132").Value = 7753.8948
$ws.Range("M132").Value = -5223.8948
$ws.Range("H141").Value = 162249.5
$ws.Range("J141").Value = 162249.5
$ws.Range("L141").Value = 162249.5
$ws.Range("N141").Value = -172609.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 399.1
$ws.Range("J12").Value = 486.5
$ws.Range("L12").Value = 1459.5
$ws.Range("N12").Value = -1805.5
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H34").Value = 61627.367
$ws.Range("J34").Value = 129944.445
$ws.Range("L34").Value = 389833.335
$ws.Range("N34").Value = -390001.335
$ws.Range("H55").Value = 9908.333000000001
$ws.Range("J55").Value = 14387.5
$ws.Range("L55").Value = 43162.5
$ws.Range("N55").Value = -43516.5
$ws.Range("H81").Value = 25457.111
$ws.Range("I81").Value = 1287.25
$ws.Range("K81").Value = 3861.75
$ws.Range("M81").Value = -2738.75
$ws.Range("H84").Value = 25457.111
$ws.Range("I84").Value = 1287.25
$ws.Range("K84").Value = 11585.25
$ws.Range("M84").Value = -5969.25
$ws.Range("H98").Value = 1244.3334
$ws.Range("I98").Value = 3090
$ws.Range("J98").Value = 960.38464
$ws.Range("K98").Value = 9270
$ws.Range("L98").Value = 2881.15392
$ws.Range("M98").Value = -7772
$ws.Range("N98").Value = -5877.15392
$ws.Range("H112").Value = 5799.6665
$ws.Range("J112").Value = 5799.5
$ws.Range("L112").Value = 17398.5
$ws.Range("N112").Value = -19614.5
$ws.Range("H128").Value = 158331
$ws.Range("I128").Value = 158331
$ws.Range("K128").Value = 474993
$ws.Range("M128").Value = -470013

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 100007816
$ws.Range("I70").Value = 6998
$ws.Range("J70").Value = 125008024
$ws.Range("K70").Value = 6998
$ws.Range("L70").Value = 125008024
$ws.Range("M70").Value = -6728
$ws.Range("N70").Value = -125008564
$ws.Range("H73").Value = 100007816
$ws.Range("I73").Value = 6998
$ws.Range("J73").Value = 125008024
$ws.Range("K73").Value = 6998
$ws.Range("L73").Value = 125008024
$ws.Range("M73").Value = -6062
$ws.Range("N73").Value = -125009896
$ws.Range("H97").Value = 17931.5
$ws.Range("I97").Value = 17931.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 17931.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -17435.5
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 1127.8
$ws.Range("I102").Value = 1127.8
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1127.8
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 494.2
$ws.Range("N102").ClearContents()
$ws.Range("H113").Value = 733660.7
$ws.Range("I113").Value = 929803.8
$ws.Range("K113").Value = 929803.8
$ws.Range("M113").Value = -927633.8
$ws.Range("H122").Value = 8011.4517
$ws.Range("I122").Value = 7104.095
$ws.Range("K122").Value = 21312.285
$ws.Range("M122").Value = -18862.285

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5603.68
$ws.Range("I7").Value = 6178.357
$ws.Range("J7").Value = 4872.273
$ws.Range("K7").Value = 6178.357
$ws.Range("L7").Value = 4872.273
$ws.Range("M7").Value = -6066.357
$ws.Range("N7").Value = -5096.273
$ws.Range("H40").Value = 3583.0857
$ws.Range("I40").Value = 3440
$ws.Range("J40").Value = 3825.2307
$ws.Range("K40").Value = 3440
$ws.Range("L40").Value = 3825.2307
$ws.Range("M40").Value = -3304
$ws.Range("N40").Value = -4097.2307
$ws.Range("H61").Value = 6647.4375
$ws.Range("I61").Value = 4636.1
$ws.Range("K61").Value = 4636.1
$ws.Range("M61").Value = -4434.1
$ws.Range("H93").Value = 1279.7646
$ws.Range("I93").Value = 1145.1538
$ws.Range("J93").Value = 1717.25
$ws.Range("K93").Value = 1145.1538
$ws.Range("L93").Value = 1717.25
$ws.Range("M93").Value = 102.8462
$ws.Range("N93").Value = -4213.25
$ws.Range("H100").Value = 3333.6
$ws.Range("I100").Value = 3333.6
$ws.Range("K100").Value = 3333.6
$ws.Range("M100").Value = -2792.6
$ws.Range("H113").Value = 6647.4375
$ws.Range("I113").Value = 4636.1
$ws.Range("K113").Value = 4636.1
$ws.Range("M113").Value = -2466.1
$ws.Range("H122").Value = 721666.5600000001
$ws.Range("J122").Value = 717890.4399999999
$ws.Range("L122").Value = 2153671.32
$ws.Range("N122").Value = -2158571.32
$ws.Range("H126").Value = 5603.68
$ws.Range("I126").Value = 6178.357
$ws.Range("J126").Value = 4872.273
$ws.Range("K126").Value = 18535.071
$ws.Range("L126").Value = 14616.819
$ws.Range("M126").Value = -16065.071
$ws.Range("N126").Value = -19556.819
$ws.Range("H132").Value = 3446.6
$ws.Range("I132").Value = 2629.6785
$ws.Range("K132").Value = 7889.0355
$ws.Range("M132").Value = -5359.0355

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 93419.82000000001
$ws.Range("J96").Value = 3243.5
$ws.Range("L96").Value = 3243.5
$ws.Range("N96").Value = -5989.5
$ws.Range("H113").Value = 1700
$ws.Range("J113").Value = 3663.3333
$ws.Range("L113").Value = 10989.9999
$ws.Range("N113").Value = -15329.9999
$ws.Range("H122").Value = 27031458
$ws.Range("I122").Value = 38465230
$ws.Range("K122").Value = 115395690
$ws.Range("M122").Value = -115393240
$ws.Range("H126").Value = 3375
$ws.Range("I126").Value = 3933.3333
$ws.Range("K126").Value = 11799.9999
$ws.Range("M126").Value = -9329.999899999999
